$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Потрачено" (spent) values for the second week's data
$ws.Range("C3").Value = 3242
$ws.Range("C4").Value = 100
$ws.Range("C5").Value = 231
$ws.Range("C6").Value = 92
$ws.Range("C7").Value = 122

# Update totals row (budget slightly revised, spent updated)
$ws.Range("B8").Value = 2999
$ws.Range("C8").Value = 501

# Move the active selection to C3
$ws.Range("C3").Select()
